# Append 17 new test-result rows (rows 46-62) to Sheet1, mirroring the
# existing "tests" table layout:
#   A = row index (0-based, bold/boxed/centered like the header)
#   B = id
#   C = user_id
#   D = test_name
#   E = language
#   F = is_finished  (text "True"/"False", same as existing rows)
#   G = result
#   H = datetime     (text "dd/mm/yyyy", same as existing rows)
#
# F/H hold strings that look like a boolean / a date, so a plain
# $cell.Value = "..." assignment gets auto-coerced by Excel's type
# sniffer. To keep them as literal text (matching the rest of the
# column) we borrow the real values from cells that already store that
# exact text: F2/F12 already hold "True"/"False", and two scratch cells
# seeded once (with a leading apostrophe) supply the two new date
# strings. Copying *values only* (xlPasteValues) carries the text type
# over without carrying any formatting, so the destination cells keep
# the table's normal (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

$newRows = @(
    @(44, 45, 3, "sogi_assessment",                 "ru", "True",  "small_result",  "12/06/2023"),
    @(45, 46, 1, "sogi_assessment",                 "ru", "True",  "small_result",  "12/06/2023"),
    @(46, 47, 3, "pkp_assessment",                  "ru", "True",  "medium_result", "12/06/2023"),
    @(47, 48, 1, "hiv_risk_assessment",              "ru", "True",  "medium_result", "12/06/2023"),
    @(48, 49, 1, "pkp_assessment",                  "ru", "True",  "medium_result", "12/06/2023"),
    @(49, 50, 3, "understanding_PLHIV_assessment",   "ru", "True",  "small_result",  "12/06/2023"),
    @(50, 51, 3, "hiv_risk_assessment",              "ru", "True",  "medium_result", "12/06/2023"),
    @(51, 52, 1, "understanding_PLHIV_assessment",   "ru", "False", "None",          "12/06/2023"),
    @(52, 53, 1, "hiv_knowledge_assessment",         "ru", "False", "None",          "12/06/2023"),
    @(53, 54, 1, "pkp_assessment",                  "ru", "True",  "medium_result", "16/06/2023"),
    @(54, 55, 1, "hiv_risk_assessment",              "ru", "True",  "small_result",  "16/06/2023"),
    @(55, 56, 1, "pkp_assessment",                  "ru", "True",  "medium_result", "16/06/2023"),
    @(56, 57, 1, "pkp_assessment",                  "ru", "True",  "small_result",  "16/06/2023"),
    @(57, 58, 1, "sogi_assessment",                 "ru", "True",  "small_result",  "16/06/2023"),
    @(58, 59, 1, "sogi_assessment",                 "ru", "True",  "small_result",  "16/06/2023"),
    @(59, 60, 1, "sogi_assessment",                 "ru", "True",  "small_result",  "16/06/2023"),
    @(60, 61, 1, "hiv_risk_assessment",              "ru", "True",  "small_result",  "16/06/2023")
)

$startRow = 46

# Existing cells already holding exactly "True" / "False" as text.
$trueSeed  = $ws.Cells.Item(2, 6)
$falseSeed = $ws.Cells.Item(12, 6)

# Seed the two new date strings (not present anywhere in the existing
# data) once each, off to the side, then fan them out with value-only
# pastes so only the scratch cells ever carry the forced-text styling;
# the scratch cells are removed again afterwards.
$scratchDec12 = $ws.Cells.Item(1, 26)   # Z1
$scratchJun16 = $ws.Cells.Item(2, 26)   # Z2
$scratchDec12.Value = "'12/06/2023"
$scratchJun16.Value = "'16/06/2023"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]

    if ($data[5] -eq "True") {
        $trueSeed.Copy()
    } else {
        $falseSeed.Copy()
    }
    $ws.Cells.Item($r, 6).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($r, 7).Value = $data[6]

    if ($data[7] -eq "12/06/2023") {
        $scratchDec12.Copy()
    } else {
        $scratchJun16.Copy()
    }
    $ws.Cells.Item($r, 8).PasteSpecial($xlPasteValues)
}

$scratchDec12.Clear()
$scratchJun16.Clear()

# Column A on every existing data row shares one bold/boxed/centered
# style (the same style used by the header row). Copy that formatting
# from the last pre-existing row (A45) onto the new A46:A62 cells so
# they match the rest of the table (re-using the existing style index
# instead of fabricating a new, merely similar-looking one).
$ws.Cells.Item(45, 1).Copy()
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0
